$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number but must remain text (matches
# the column's existing text formatting, e.g. '1.00', '0.501') get their
# number format forced to Text ('@') before the value is written, exactly
# as Excel would require to keep a numeric-looking entry as a string.
$textCells = @("D5","D6","D9","D10","D11","D12","D15","D17","D20","D21","D22","D26","D28","D29","D32","D34","D36","D37","D38","D43","D44","D47","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price, 1h volume, and the handful of
# rows whose rank swapped with a neighbor) cell by cell.
$ws.Range("D2").Value = '25.705.58'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '1.626.42'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '214.19'
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("D6").Value = '0.501'
$ws.Range("E6").Value = '  -1.16%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("D9").Value = '0.0635'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("D10").Value = '19.50'
$ws.Range("E10").Value = '  -4.84%  '
$ws.Range("D11").Value = '0.0781'
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '4.23'
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.624.94'
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("D14").Value = '1.851.23'
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").Value = '0.551'
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("D16").Value = '0.0₃0761'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").Value = '62.58'
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("D18").Value = '25.727.82'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").Value = '4.41'
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").Value = '193.19'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").Value = '9.91'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").Value = '139.40'
$ws.Range("E26").Value = '  -1.35%  '
$ws.Range("D28").Value = '6.83'
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").Value = '15.44'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("D32").Value = '3.31'
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").Value = '1.58'
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").Value = '0.895'
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").Value = '2.56'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").Value = '0.545'
$ws.Range("E38").Value = '  -1.91%  '
$ws.Range("D39").Value = '1.108.41'
$ws.Range("E39").Value = '  -2.19%  '
$ws.Range("E40").Value = '  -1.20%  '
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("D43").Value = '100.02'
$ws.Range("E43").Value = '  +0.93%  '
$ws.Range("D44").Value = '0.797'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").Value = '1.758.10'
$ws.Range("E45").Value = '  -1.31%  '
$ws.Range("D46").Value = '0.0₆0106'
$ws.Range("E46").Value = '  -4.35%  '
$ws.Range("D47").Value = '54.82'
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '7.71'
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '0.418'
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("E50").Value = '  +3.23%  '
$ws.Range("E51").Value = '  -0.73%  '
